$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.120.11"
$ws.Range("E2").Value = "  +5.73%  "
$ws.Range("D3").Value = "2.264.83"
$ws.Range("E3").Value = "  +2.53%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E5").Value = "  +3.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.18"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +7.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.534"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +4.63%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  +4.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.92"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +8.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.84"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +9.66%  "
$ws.Range("E12").Value = "  +3.14%  "
$ws.Range("E13").Value = "  +3.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.71"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.63%  "
$ws.Range("D15").Value = "2.615.57"
$ws.Range("E15").Value = "  +2.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.18"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.72%  "
$ws.Range("D17").Value = "2.270.03"
$ws.Range("E17").Value = "  +4.38%  "
$ws.Range("E18").Value = "  +4.24%  "
$ws.Range("D19").Value = "42.010.02"
$ws.Range("E19").Value = "  +5.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.24"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +9.47%  "
$ws.Range("E21").Value = "  +3.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.96"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.49"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "242.69"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.59"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.85%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +5.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.01"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.19"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.23%  "
$ws.Range("E30").Value = "  +5.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.32"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +8.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "158.50"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  +5.27%  "
$ws.Range("E35").Value = "  +5.16%  "
$ws.Range("E36").Value = "  +6.50%  "
$ws.Range("E37").Value = "  +3.22%  "
$ws.Range("E38").Value = "  +7.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.64"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +9.81%  "
$ws.Range("E40").Value = "  +4.38%  "
$ws.Range("E41").Value = "  +6.13%  "
$ws.Range("E42").Value = "  +7.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.13"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +13.83%  "
$ws.Range("D44").Value = "2.056.17"
$ws.Range("E44").Value = "  -2.51%  "
$ws.Range("E45").Value = "  +5.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.11"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.39%  "
$ws.Range("E47").Value = "  +9.13%  "
$ws.Range("E48").Value = "  -4.78%  "
$ws.Range("D49").Value = "2.488.26"
$ws.Range("E49").Value = "  +2.33%  "
$ws.Range("E50").Value = "  +3.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.09"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +6.84%  "
